$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27: num_customers 53 -> 54, retention_rate recalculated as C27/D27
$ws.Range("C27").Value = 54
$ws.Range("E27").Value = 54/2252

# Row 34: num_customers 87 -> 88, retention_rate recalculated as C34/D34
$ws.Range("C34").Value = 88
$ws.Range("E34").Value = 88/2256

# Row 37: num_customers 941 -> 945, cohort_size 941 -> 945 (retention_rate stays 1)
$ws.Range("C37").Value = 945
$ws.Range("D37").Value = 945
